# Generate Report for Handoff
# Update the "Latest Handoff" timestamps for the 9e743de5-aa8f-4758-a1b9-9c20d02a02f4
# entry (row 7) across the Overview, zh-cn and de-de sheets, reflecting a fresh handoff.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("D7").Value = "2016-03-24 11:31:58"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E7").Value = "2016-03-24 11:31:53"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E7").Value = "2016-03-24 11:31:58"
